$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring row 24 formatting in line with the other data rows (date style on A,
# percent style on E) by copying the row above before overwriting values.
$ws.Range("A23:E23").Copy($ws.Range("A24:E24"))

# Preliminary Boston 2024 (#1086) attendance numbers.
$ws.Range("A24").Value = 45570
$ws.Range("B24").Value = "SQL Saturday Boston 2024 (1086)"
$ws.Range("C24").Value = 250
$ws.Range("D24").Value = 203
$ws.Range("E24").Formula = "=IF(C24=0,0,+(C24-D24)/C24)"

# Selection moved on to D25 when the edit was saved.
[void]$ws.Range("D25").Select()
